# Update cryptos list: refresh prices and 1h volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.683.02'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '3.380.85'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '573.99'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.31'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '3.380.25'
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.45'
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').Value = '3.955.93'
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('E15').Value = '  -1.85%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '25.93'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('D17').Value = '3.381.90'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '61.805.73'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.92'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.94'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '376.87'
$ws.Range('E22').Value = '  -2.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.556'
$ws.Range('E23').Value = '  -2.69%  '
$ws.Range('D24').Value = '3.523.68'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '71.29'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('E28').Value = '  +9.96%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.61'
$ws.Range('E29').Value = '  -3.28%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +4.10%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.22'
$ws.Range('E32').Value = '  -0.87%  '
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.20'
$ws.Range('E36').Value = '  -6.15%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.54'
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.83'
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '165.02'
$ws.Range('E39').Value = '  +2.48%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0769'
$ws.Range('E40').Value = '  -3.18%  '
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '41.58'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.38'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '24.31'
$ws.Range('E47').Value = '  +3.58%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '6.86'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '23.06'
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('D50').Value = '2.384.68'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0264'
$ws.Range('E51').Value = '  -0.02%  '
